# Update project part 2
# Adds a "Non-Fishers" column (E) computed as Total Labour Force - Fishing Occupations,
# and updates the selection to the new column's data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column E
$ws.Range("E1").Value = "Non-Fishers"

# E2 gets its own (non-shared) formula
$ws.Range("E2").Formula = "=C2-B2"

# E3:E10 share one formula group (mirrors Excel's fill-down behaviour)
$ws.Range("E3:E10").Formula = "=C3-B3"

# Update selection to mirror the diff (active cell E2, selection E2:E10)
$ws.Range("E2:E10").Select()
